# ProyectoPIA5B3.docx edit script
# - "Recursos:" row: merge the two bold runs into one, center the content
#   cell vertically and fill it in with the Pictoblox bullet notes.
# - "Proceso:" row: merge the two bold runs into one, center the content
#   cell vertically and fill it in with the 14-step numbered process list
#   (a brand new numbering definition, same as Word allocates automatically
#   the first time ApplyNumberDefault() is used on a fresh paragraph).
# - "Discusión y conclusiones del proyecto:" row: merge the two bold runs
#   into one and mark a lastRenderedPageBreak on it.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Row 11 - "Recursos:"
# ---------------------------------------------------------------------

$recCell1 = $t.Cell(11, 1)
$recCell1.Range.InsertXML(@"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>
<w:p w14:paraId="15C9B5D2" w14:textId="04311AAA" w:rsidR="00F73823" w:rsidRPr="00402889" w:rsidRDefault="00AB69BF">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Recursos:</w:t>
  </w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@)

$t = $d.Tables.Item(1)
$recCell2 = $t.Cell(11, 2)
$recCell2.VerticalAlignment = 1
$recCell2.Range.InsertXML(@"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">&#8226; </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Pictoblox</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>&#8226; Utilización de algún fondo para hacer mucho más entretenido el ingreso de datos.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>&#8226; Utilización más de herramientas propias de la aplicación para simular de que nos quede mucho mejor visualmente.</w:t>
  </w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@)

# ---------------------------------------------------------------------
# Row 12 - "Proceso:"
# ---------------------------------------------------------------------

$t = $d.Tables.Item(1)
$procCell1 = $t.Cell(12, 1)
$procCell1.Range.InsertXML(@"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>
<w:p w14:paraId="62E24A9E" w14:textId="6AD8944B" w:rsidR="00F73823" w:rsidRPr="00402889" w:rsidRDefault="005D0865">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Proceso:</w:t>
  </w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@)

# Allocate a brand new list-numbering definition (this is what Word itself
# does under the hood the first time a fresh/default numbered list gets
# applied) so the content below can reference it as numId 6.
$t = $d.Tables.Item(1)
$procCell2 = $t.Cell(12, 2)
$procPara = $procCell2.Range.Paragraphs.Item(1)
$procPara.Range.ListFormat.ApplyNumberDefault()

$t = $d.Tables.Item(1)
$procCell2 = $t.Cell(12, 2)
$procCell2.VerticalAlignment = 1
$procCell2.Range.InsertXML(@"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Contextualizar el objetivo del proyecto.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Abrir PictoBlox</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Creamos un nuevo proyecto</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Agregamos el plugin de inteligencia artificial el cual nos va a ayudar hacer el reconocimiento de voz.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Agrega el control al hacer clic</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Crear la variable con la que va a trabajar edad</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Guardar estos números en una lista nos seria de mucha utilizada para más adelante poder recorrer y obtener la frecuencia absoluta.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Creamos dos variables como por ejemplo edades y frecuencia.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Utilizamos el bloque de repetir por siempre.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Usamos el bloque de pregunta, para preguntar las edades de los participantes luego lo guardamos en una lista creada anteriormente con el nombre de edades.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Para finalizar el proceso de ingreso de información se puede configurar el bloque que se detenga cuando se presiona alguna tecla y se empezaría con el nuevo proceso que es de recorrer la lista</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t>Para ello utilizamos el bloque repetir hasta que se llegue al final de la lista e ir sacando los datos en variables por ejemplo si la edad de 10 años se repite podemos ir guardando en número de repeticiones.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">Una vez terminado el conteo de datos presionamos alguna tecla que es la que nos va a ayudar con el grafico. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:ind w:left="275" w:hanging="283"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:lang w:eastAsia="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">Para poder realizar los gráficos es necesario crear clones de los objetos que se vayan a utilizar para representar y así simular como crecen los cuadros estadísticos. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@)

# ---------------------------------------------------------------------
# Row 13 - "Discusión y conclusiones del proyecto:"
# ---------------------------------------------------------------------

$t = $d.Tables.Item(1)
$discCell1 = $t.Cell(13, 1)
$discCell1.Range.InsertXML(@"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>
<w:p w14:paraId="71610C7E" w14:textId="11BDE55A" w:rsidR="00F73823" w:rsidRPr="00402889" w:rsidRDefault="002B42D5">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Discusión y conclusiones del proyecto:</w:t>
  </w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@)

Write-Output "Edit complete"
